$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegistrationData")

# Clear the phone number value in D6 but keep its existing formatting
$ws.Range("D6").ClearContents()

# Add hyperlinks for the new row 7 entries first
$ws.Hyperlinks.Add($ws.Range("E7"), "mailto:Gemini@2515") | Out-Null
$ws.Hyperlinks.Add($ws.Range("F7"), "mailto:Gemini@251522") | Out-Null

# Copy formatting from row 6 into row 7 (after hyperlink creation so it isn't overwritten)
$ws.Range("E6:G6").Copy()
$ws.Range("E7").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = 0

# Re-apply the values/text for row 7
$ws.Range("E7").Value = "Gemini@2515"
$ws.Range("F7").Value = "Gemini@251522"
$ws.Range("G7").Value = "Invalid Data"
$ws.Range("H7").Value = "Need to provide data"

# Update the active selection
$ws.Range("D6").Select() | Out-Null
